$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 131025505
$ws.Range("B2").Value = 79243
$ws.Range("E2").Value = 6425
$ws.Range("F2").Value = 'Garnlav'
$ws.Range("G2").Value = 'Alectoria sarmentosa'
$ws.Range("H2").Value = '(Ach.) Ach.'
$ws.Range("P2").Value = 'Sjöarsjön Nordöst, Jmt'
$ws.Range("Q2").Value = 464214
$ws.Range("R2").Value = 7042578
$ws.Range("AC2").Value = 'På flera granar.'
$ws.Range("AW2").Value = 'Kristian Zackrisson'
$ws.Range("AX2").Value = 'Kristian Zackrisson'

# Row 3
$ws.Range("A3").Value = 131025514
$ws.Range("Q3").Value = 464434
$ws.Range("R3").Value = 7042377
$ws.Range("AC3").Value = 'På flera granar vid en lucka.'

# Row 4
$ws.Range("A4").Value = 131024454
$ws.Range("B4").Value = 57884
$ws.Range("E4").Value = 100109
$ws.Range("F4").Value = 'Tretåig hackspett'
$ws.Range("G4").Value = 'Picoides tridactylus'
$ws.Range("H4").Value = '(Linnaeus, 1758)'
$ws.Range("P4").Value = 'Sjöarsjön N, Jmt'
$ws.Range("Q4").Value = 464321
$ws.Range("R4").Value = 7042482
$ws.Range("AC4").Value = 'Ringhack äldre'
$ws.Range("AW4").Value = 'Benny Öwre'
$ws.Range("AX4").Value = 'Benny Öwre'

# Row 18
$ws.Range("A18").Value = 131024468
$ws.Range("B18").Value = 57884
$ws.Range("E18").Value = 100109
$ws.Range("F18").Value = 'Tretåig hackspett'
$ws.Range("G18").Value = 'Picoides tridactylus'
$ws.Range("H18").Value = '(Linnaeus, 1758)'
$ws.Range("P18").Value = 'Sjöarsjön N, Jmt'
$ws.Range("Q18").Value = 464332
$ws.Range("R18").Value = 7042344
$ws.Range("AC18").Value = 'Ringhack färska och äldre'
$ws.Range("AW18").Value = 'Benny Öwre'
$ws.Range("AX18").Value = 'Benny Öwre'

# Row 19
$ws.Range("A19").Value = 131024458
$ws.Range("B19").Value = 57884
$ws.Range("E19").Value = 100109
$ws.Range("F19").Value = 'Tretåig hackspett'
$ws.Range("G19").Value = 'Picoides tridactylus'
$ws.Range("H19").Value = '(Linnaeus, 1758)'
$ws.Range("P19").Value = 'Sjöarsjön N, Jmt'
$ws.Range("Q19").Value = 464405
$ws.Range("R19").Value = 7042391
$ws.Range("AC19").Value = 'Ringhack'
$ws.Range("AW19").Value = 'Benny Öwre'
$ws.Range("AX19").Value = 'Benny Öwre'

# Row 20
$ws.Range("A20").Value = 131024464
$ws.Range("Q20").Value = 464337
$ws.Range("R20").Value = 7042380

# Row 21
$ws.Range("A21").Value = 131025468
$ws.Range("M21").Value = 'äldre spår'
$ws.Range("P21").Value = 'Sjöarsjön Nordöst, Jmt'
$ws.Range("Q21").Value = 464445
$ws.Range("R21").Value = 7042578
$ws.Range("AC21").Value = 'Ringhack, äldre, på en gran.'
$ws.Range("AH21").Value = 'Granskog'
$ws.Range("AJ21").Value = 'gran'
$ws.Range("AK21").Value = 'Picea abies'
$ws.Range("AM21").Value = 'Trädstam på levande träd'
$ws.Range("AO21").Value = 'Stem on living tree # Picea abies'
$ws.Range("AW21").Value = 'Kristian Zackrisson'
$ws.Range("AX21").Value = 'Kristian Zackrisson'

# Row 22
$ws.Range("A22").Value = 131025510
$ws.Range("B22").Value = 79243
$ws.Range("E22").Value = 6425
$ws.Range("F22").Value = 'Garnlav'
$ws.Range("G22").Value = 'Alectoria sarmentosa'
$ws.Range("H22").Value = '(Ach.) Ach.'
$ws.Range("P22").Value = 'Sjöarsjön Nordöst, Jmt'
$ws.Range("Q22").Value = 464463
$ws.Range("R22").Value = 7042542
$ws.Range("AC22").Value = 'På flera granar.'
$ws.Range("AW22").Value = 'Kristian Zackrisson'
$ws.Range("AX22").Value = 'Kristian Zackrisson'

# Row 23
$ws.Range("A23").Value = 131025525
$ws.Range("B23").Value = 79243
$ws.Range("E23").Value = 6425
$ws.Range("F23").Value = 'Garnlav'
$ws.Range("G23").Value = 'Alectoria sarmentosa'
$ws.Range("H23").Value = '(Ach.) Ach.'
$ws.Range("M23").Value = $null
$ws.Range("Q23").Value = 464067
$ws.Range("R23").Value = 7042585
$ws.Range("AC23").Value = 'På gran.'
$ws.Range("AH23").Value = $null
$ws.Range("AJ23").Value = $null
$ws.Range("AK23").Value = $null
$ws.Range("AM23").Value = $null
$ws.Range("AO23").Value = $null

# Row 55
$ws.Range("A55").Value = 131024463
$ws.Range("B55").Value = 57884
$ws.Range("E55").Value = 100109
$ws.Range("F55").Value = 'Tretåig hackspett'
$ws.Range("G55").Value = 'Picoides tridactylus'
$ws.Range("H55").Value = '(Linnaeus, 1758)'
$ws.Range("Q55").Value = 464371
$ws.Range("R55").Value = 7042397
$ws.Range("AC55").Value = 'Ringhack äldre'

# Row 56
$ws.Range("A56").Value = 131025483
$ws.Range("B56").Value = 57884
$ws.Range("E56").Value = 100109
$ws.Range("F56").Value = 'Tretåig hackspett'
$ws.Range("G56").Value = 'Picoides tridactylus'
$ws.Range("H56").Value = '(Linnaeus, 1758)'
$ws.Range("M56").Value = 'färska spår'
$ws.Range("P56").Value = 'Sjöarsjön Nordöst, Jmt'
$ws.Range("Q56").Value = 464021
$ws.Range("R56").Value = 7042325
$ws.Range("AC56").Value = 'Ringhack, färska, längs minst 6 meter på en gran ca 4 meter från strandkanten. Rikligt med spår av rinnande kåda/sav på granens stam.'
$ws.Range("AH56").Value = 'Granskog'
$ws.Range("AJ56").Value = 'gran'
$ws.Range("AK56").Value = 'Picea abies'
$ws.Range("AM56").Value = 'Trädstam på levande träd'
$ws.Range("AO56").Value = 'Stem on living tree # Picea abies'
$ws.Range("AW56").Value = 'Kristian Zackrisson'
$ws.Range("AX56").Value = 'Kristian Zackrisson'

# Row 57
$ws.Range("A57").Value = 131025469
$ws.Range("B57").Value = 57884
$ws.Range("E57").Value = 100109
$ws.Range("F57").Value = 'Tretåig hackspett'
$ws.Range("G57").Value = 'Picoides tridactylus'
$ws.Range("H57").Value = '(Linnaeus, 1758)'
$ws.Range("M57").Value = 'äldre spår'
$ws.Range("P57").Value = 'Sjöarsjön Nordöst, Jmt'
$ws.Range("Q57").Value = 464464
$ws.Range("R57").Value = 7042544
$ws.Range("AC57").Value = 'Ringhack, äldre, på en granhögstubbe.'
$ws.Range("AH57").Value = 'Granskog'
$ws.Range("AJ57").Value = 'gran'
$ws.Range("AK57").Value = 'Picea abies'
$ws.Range("AM57").Value = 'Stående död trädstam/högstubbe'
$ws.Range("AO57").Value = 'Standing dead tree/snags # Picea abies'
$ws.Range("AW57").Value = 'Kristian Zackrisson'
$ws.Range("AX57").Value = 'Kristian Zackrisson'

# Row 58
$ws.Range("A58").Value = 131024511
$ws.Range("B58").Value = 91828
$ws.Range("E58").Value = 5432
$ws.Range("F58").Value = 'Granticka'
$ws.Range("G58").Value = 'Porodaedalea chrysoloma s.lat.'
$ws.Range("H58").Value = $null
$ws.Range("P58").Value = 'Sjöarsjön N, Jmt'
$ws.Range("Q58").Value = 463866
$ws.Range("R58").Value = 7042468
$ws.Range("AC58").Value = $null
$ws.Range("AW58").Value = 'Benny Öwre'
$ws.Range("AX58").Value = 'Benny Öwre'

# Row 59
$ws.Range("A59").Value = 131024510
$ws.Range("B59").Value = 91828
$ws.Range("E59").Value = 5432
$ws.Range("F59").Value = 'Granticka'
$ws.Range("G59").Value = 'Porodaedalea chrysoloma s.lat.'
$ws.Range("H59").Value = $null
$ws.Range("Q59").Value = 464031
$ws.Range("R59").Value = 7042403
$ws.Range("AC59").Value = $null

# Row 60
$ws.Range("A60").Value = 131024507
$ws.Range("B60").Value = 91828
$ws.Range("E60").Value = 5432
$ws.Range("F60").Value = 'Granticka'
$ws.Range("G60").Value = 'Porodaedalea chrysoloma s.lat.'
$ws.Range("H60").Value = $null
$ws.Range("M60").Value = $null
$ws.Range("P60").Value = 'Sjöarsjön N, Jmt'
$ws.Range("Q60").Value = 464395
$ws.Range("R60").Value = 7042406
$ws.Range("AC60").Value = $null
$ws.Range("AH60").Value = $null
$ws.Range("AJ60").Value = $null
$ws.Range("AK60").Value = $null
$ws.Range("AM60").Value = $null
$ws.Range("AO60").Value = $null
$ws.Range("AW60").Value = 'Benny Öwre'
$ws.Range("AX60").Value = 'Benny Öwre'

# Row 61
$ws.Range("A61").Value = 131025504
$ws.Range("B61").Value = 79243
$ws.Range("E61").Value = 6425
$ws.Range("F61").Value = 'Garnlav'
$ws.Range("G61").Value = 'Alectoria sarmentosa'
$ws.Range("H61").Value = '(Ach.) Ach.'
$ws.Range("M61").Value = $null
$ws.Range("Q61").Value = 464204
$ws.Range("R61").Value = 7042559
$ws.Range("AC61").Value = 'På gran.'
$ws.Range("AH61").Value = $null
$ws.Range("AJ61").Value = $null
$ws.Range("AK61").Value = $null
$ws.Range("AM61").Value = $null
$ws.Range("AO61").Value = $null

# Row 73
$ws.Range("A73").Value = 131024453
$ws.Range("B73").Value = 57884
$ws.Range("E73").Value = 100109
$ws.Range("F73").Value = 'Tretåig hackspett'
$ws.Range("G73").Value = 'Picoides tridactylus'
$ws.Range("H73").Value = '(Linnaeus, 1758)'
$ws.Range("P73").Value = 'Sjöarsjön N, Jmt'
$ws.Range("Q73").Value = 464327
$ws.Range("R73").Value = 7042477
$ws.Range("AC73").Value = 'Ringhack färska och äldre'
$ws.Range("AW73").Value = 'Benny Öwre'
$ws.Range("AX73").Value = 'Benny Öwre'

# Row 74
$ws.Range("A74").Value = 131024451
$ws.Range("Q74").Value = 464275
$ws.Range("R74").Value = 7042455

# Row 75
$ws.Range("A75").Value = 131025480
$ws.Range("M75").Value = 'äldre spår'
$ws.Range("P75").Value = 'Sjöarsjön Nordöst, Jmt'
$ws.Range("Q75").Value = 464378
$ws.Range("R75").Value = 7042302
$ws.Range("AC75").Value = 'Ringhack, äldre, på stambasen av en stående död gran (BHD >10cm).'
$ws.Range("AH75").Value = 'Granskog'
$ws.Range("AJ75").Value = 'gran'
$ws.Range("AK75").Value = 'Picea abies'
$ws.Range("AM75").Value = 'Stående död trädstam/högstubbe'
$ws.Range("AO75").Value = 'Standing dead tree/snags # Picea abies'
$ws.Range("AW75").Value = 'Kristian Zackrisson'
$ws.Range("AX75").Value = 'Kristian Zackrisson'

# Row 76
$ws.Range("A76").Value = 131025521
$ws.Range("B76").Value = 79243
$ws.Range("E76").Value = 6425
$ws.Range("F76").Value = 'Garnlav'
$ws.Range("G76").Value = 'Alectoria sarmentosa'
$ws.Range("H76").Value = '(Ach.) Ach.'
$ws.Range("M76").Value = $null
$ws.Range("Q76").Value = 464007
$ws.Range("R76").Value = 7042409
$ws.Range("AC76").Value = 'Långväxta bålar på gran.'
$ws.Range("AH76").Value = $null
$ws.Range("AJ76").Value = $null
$ws.Range("AK76").Value = $null
$ws.Range("AM76").Value = $null
$ws.Range("AO76").Value = $null

# Row 94
$ws.Range("A94").Value = 131024512
$ws.Range("B94").Value = 91828
$ws.Range("E94").Value = 5432
$ws.Range("F94").Value = 'Granticka'
$ws.Range("G94").Value = 'Porodaedalea chrysoloma s.lat.'
$ws.Range("H94").Value = $null
$ws.Range("M94").Value = $null
$ws.Range("P94").Value = 'Sjöarsjön N, Jmt'
$ws.Range("Q94").Value = 464059
$ws.Range("R94").Value = 7042514
$ws.Range("AC94").Value = $null
$ws.Range("AH94").Value = $null
$ws.Range("AI94").Value = $null
$ws.Range("AJ94").Value = $null
$ws.Range("AK94").Value = $null
$ws.Range("AM94").Value = $null
$ws.Range("AO94").Value = $null
$ws.Range("AW94").Value = 'Benny Öwre'
$ws.Range("AX94").Value = 'Benny Öwre'

# Row 95
$ws.Range("A95").Value = 131025473
$ws.Range("B95").Value = 57884
$ws.Range("E95").Value = 100109
$ws.Range("F95").Value = 'Tretåig hackspett'
$ws.Range("G95").Value = 'Picoides tridactylus'
$ws.Range("H95").Value = '(Linnaeus, 1758)'
$ws.Range("M95").Value = 'färska spår'
$ws.Range("P95").Value = 'Sjöarsjön Nordöst, Jmt'
$ws.Range("Q95").Value = 464413
$ws.Range("R95").Value = 7042312
$ws.Range("AC95").Value = 'Ringhack, färska och äldre, längs ett par meter på en gran med spår av rikligt kåda-flöde.'
$ws.Range("AH95").Value = 'Granskog'
$ws.Range("AI95").Value = 'Gammal granskog.'
$ws.Range("AJ95").Value = 'gran'
$ws.Range("AK95").Value = 'Picea abies'
$ws.Range("AM95").Value = 'Trädstam på levande träd'
$ws.Range("AO95").Value = 'Stem on living tree # Picea abies'
$ws.Range("AW95").Value = 'Kristian Zackrisson'
$ws.Range("AX95").Value = 'Kristian Zackrisson'
